$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: the "Doi model moi: " bullet item loses its numbered-list
# formatting (numPr ilvl=0/numId=2) in favour of a plain indent, and its
# text is replaced with the new instructions.
# ---------------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*model m*i*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the 'Doi model moi' paragraph"
}

$xmlReplacement = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:ind w:left="1080"/>
<w:rPr>
<w:lang w:val="vi-VN"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:lang w:val="vi-VN"/>
</w:rPr>
<w:t xml:space="preserve">Tạo một project mới và tạo API key, sau đó: </w:t>
</w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$target.Range.InsertXML($xmlReplacement)

# ---------------------------------------------------------------------------
# Change 2: a brand-new bulleted paragraph ("streamlit run app/ui_streamlit.py")
# is appended right after the paragraph that ends with the closing quote
# ( .\env\set_gemini.ps1 ... "gemini-2.5-flash" ), and right before the
# trailing empty paragraph that closes the document.
# ---------------------------------------------------------------------------
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*set_gemini.ps1*") {
        $anchor = $p
    }
}

if ($anchor -eq $null) {
    throw "Could not locate the 'set_gemini.ps1' paragraph"
}

# Insert point: right before the paragraph mark that ends $anchor, so the
# new paragraph is created between $anchor and whatever already follows it
# (the pre-existing empty trailing paragraph), instead of merging into it.
$insertPos = $anchor.Range.End - 1
$insertRange = $d.Range($insertPos, $insertPos)

$newParaXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr>
<w:ilvl w:val="0"/>
<w:numId w:val="2"/>
</w:numPr>
<w:rPr>
<w:lang w:val="vi-VN"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:lang w:val="vi-VN"/>
</w:rPr>
<w:t>streamlit run app/ui_streamlit.py</w:t>
</w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$insertRange.InsertXML($newParaXml)
